$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clean up the header labels -------------------------------------
# Row 1 (columns B..DQ) and Column A (rows 2..121) currently hold Python
# tuple-repr strings like "('ANG_1',)". Strip the wrapper so they just
# read "ANG_1".
for ($c = 2; $c -le 121; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $txt = $cell.Value2
    if ($txt -ne $null -and $txt -match "^\('(.*)',\)$") {
        $cell.Value = $matches[1]
    }
}

for ($r = 2; $r -le 121; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $txt = $cell.Value2
    if ($txt -ne $null -and $txt -match "^\('(.*)',\)$") {
        $cell.Value = $matches[1]
    }
}

# --- 2. Zero out every non-zero value in the correlation matrix --------
for ($r = 2; $r -le 121; $r++) {
    for ($c = 2; $c -le 121; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne 0) {
            $cell.Value = 0
        }
    }
}
